$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A117:A130").NumberFormat = "@"
$ws.Range("A117").Value = "2026-01-28"
$ws.Range("B117").Value = "16:46:48"
$ws.Range("C117").Value = "16:00"
$ws.Range("D117").Value = "Bathroom"
$ws.Range("E117").Value = "No Motion"
$ws.Range("F117").Value = "Inactive"
$ws.Range("A118").Value = "2026-01-28"
$ws.Range("B118").Value = "16:46:49"
$ws.Range("C118").Value = "16:00"
$ws.Range("D118").Value = "Bathroom"
$ws.Range("E118").Value = "No Motion"
$ws.Range("F118").Value = "Inactive"
$ws.Range("A119").Value = "2026-01-28"
$ws.Range("B119").Value = "16:46:51"
$ws.Range("C119").Value = "16:00"
$ws.Range("D119").Value = "Bathroom"
$ws.Range("E119").Value = "No Motion"
$ws.Range("F119").Value = "Inactive"
$ws.Range("A120").Value = "2026-01-28"
$ws.Range("B120").Value = "16:46:56"
$ws.Range("C120").Value = "16:00"
$ws.Range("D120").Value = "Bathroom"
$ws.Range("E120").Value = "No Motion"
$ws.Range("F120").Value = "Inactive"
$ws.Range("A121").Value = "2026-01-28"
$ws.Range("B121").Value = "16:47:02"
$ws.Range("C121").Value = "16:00"
$ws.Range("D121").Value = "Bathroom"
$ws.Range("E121").Value = "No Motion"
$ws.Range("F121").Value = "Inactive"
$ws.Range("A122").Value = "2026-01-28"
$ws.Range("B122").Value = "16:47:07"
$ws.Range("C122").Value = "16:00"
$ws.Range("D122").Value = "Bathroom"
$ws.Range("E122").Value = "No Motion"
$ws.Range("F122").Value = "Inactive"
$ws.Range("A123").Value = "2026-01-28"
$ws.Range("B123").Value = "16:47:12"
$ws.Range("C123").Value = "16:00"
$ws.Range("D123").Value = "Bathroom"
$ws.Range("E123").Value = "No Motion"
$ws.Range("F123").Value = "Inactive"
$ws.Range("A124").Value = "2026-01-28"
$ws.Range("B124").Value = "16:47:17"
$ws.Range("C124").Value = "16:00"
$ws.Range("D124").Value = "Bathroom"
$ws.Range("E124").Value = "No Motion"
$ws.Range("F124").Value = "Inactive"
$ws.Range("A125").Value = "2026-01-28"
$ws.Range("B125").Value = "16:47:22"
$ws.Range("C125").Value = "16:00"
$ws.Range("D125").Value = "Bathroom"
$ws.Range("E125").Value = "No Motion"
$ws.Range("F125").Value = "Inactive"
$ws.Range("A126").Value = "2026-01-28"
$ws.Range("B126").Value = "16:47:27"
$ws.Range("C126").Value = "16:00"
$ws.Range("D126").Value = "Bathroom"
$ws.Range("E126").Value = "No Motion"
$ws.Range("F126").Value = "Inactive"
$ws.Range("A127").Value = "2026-01-28"
$ws.Range("B127").Value = "16:47:32"
$ws.Range("C127").Value = "16:00"
$ws.Range("D127").Value = "Bathroom"
$ws.Range("E127").Value = "No Motion"
$ws.Range("F127").Value = "Inactive"
$ws.Range("A128").Value = "2026-01-28"
$ws.Range("B128").Value = "16:47:37"
$ws.Range("C128").Value = "16:00"
$ws.Range("D128").Value = "Bathroom"
$ws.Range("E128").Value = "No Motion"
$ws.Range("F128").Value = "Inactive"
$ws.Range("A129").Value = "2026-01-28"
$ws.Range("B129").Value = "16:47:42"
$ws.Range("C129").Value = "16:00"
$ws.Range("D129").Value = "Bathroom"
$ws.Range("E129").Value = "No Motion"
$ws.Range("F129").Value = "Inactive"
$ws.Range("A130").Value = "2026-01-28"
$ws.Range("B130").Value = "16:47:47"
$ws.Range("C130").Value = "16:00"
$ws.Range("D130").Value = "Bathroom"
$ws.Range("E130").Value = "No Motion"
$ws.Range("F130").Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A117:A127").NumberFormat = "@"
$ws.Range("E117:E127").NumberFormat = "@"
$ws.Range("A117").Value = "2026-01-28"
$ws.Range("B117").Value = "16:46:49"
$ws.Range("C117").Value = "16:00"
$ws.Range("D117").Value = "Bathroom"
$ws.Range("E117").Value = "86.7%"
$ws.Range("F117").Value = "Active"
$ws.Range("A118").Value = "2026-01-28"
$ws.Range("B118").Value = "16:46:50"
$ws.Range("C118").Value = "16:00"
$ws.Range("D118").Value = "Bathroom"
$ws.Range("E118").Value = "87.6%"
$ws.Range("F118").Value = "Active"
$ws.Range("A119").Value = "2026-01-28"
$ws.Range("B119").Value = "16:46:53"
$ws.Range("C119").Value = "16:00"
$ws.Range("D119").Value = "Bathroom"
$ws.Range("E119").Value = "87.6%"
$ws.Range("F119").Value = "Active"
$ws.Range("A120").Value = "2026-01-28"
$ws.Range("B120").Value = "16:47:01"
$ws.Range("C120").Value = "16:00"
$ws.Range("D120").Value = "Bathroom"
$ws.Range("E120").Value = "87.6%"
$ws.Range("F120").Value = "Active"
$ws.Range("A121").Value = "2026-01-28"
$ws.Range("B121").Value = "16:47:05"
$ws.Range("C121").Value = "16:00"
$ws.Range("D121").Value = "Bathroom"
$ws.Range("E121").Value = "86.8%"
$ws.Range("F121").Value = "Active"
$ws.Range("A122").Value = "2026-01-28"
$ws.Range("B122").Value = "16:47:09"
$ws.Range("C122").Value = "16:00"
$ws.Range("D122").Value = "Bathroom"
$ws.Range("E122").Value = "87.7%"
$ws.Range("F122").Value = "Active"
$ws.Range("A123").Value = "2026-01-28"
$ws.Range("B123").Value = "16:47:13"
$ws.Range("C123").Value = "16:00"
$ws.Range("D123").Value = "Bathroom"
$ws.Range("E123").Value = "87.7%"
$ws.Range("F123").Value = "Active"
$ws.Range("A124").Value = "2026-01-28"
$ws.Range("B124").Value = "16:47:29"
$ws.Range("C124").Value = "16:00"
$ws.Range("D124").Value = "Bathroom"
$ws.Range("E124").Value = "87.7%"
$ws.Range("F124").Value = "Active"
$ws.Range("A125").Value = "2026-01-28"
$ws.Range("B125").Value = "16:47:33"
$ws.Range("C125").Value = "16:00"
$ws.Range("D125").Value = "Bathroom"
$ws.Range("E125").Value = "87.7%"
$ws.Range("F125").Value = "Active"
$ws.Range("A126").Value = "2026-01-28"
$ws.Range("B126").Value = "16:47:37"
$ws.Range("C126").Value = "16:00"
$ws.Range("D126").Value = "Bathroom"
$ws.Range("E126").Value = "86.8%"
$ws.Range("F126").Value = "Active"
$ws.Range("A127").Value = "2026-01-28"
$ws.Range("B127").Value = "16:47:41"
$ws.Range("C127").Value = "16:00"
$ws.Range("D127").Value = "Bathroom"
$ws.Range("E127").Value = "87.7%"
$ws.Range("F127").Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A116:A127").NumberFormat = "@"
$ws.Range("A116").Value = "2026-01-28"
$ws.Range("B116").Value = "16:46:48"
$ws.Range("C116").Value = "16:00"
$ws.Range("D116").Value = "Bathroom"
$ws.Range("E116").Value = "22.8C"
$ws.Range("F116").Value = "Active"
$ws.Range("A117").Value = "2026-01-28"
$ws.Range("B117").Value = "16:46:49"
$ws.Range("C117").Value = "16:00"
$ws.Range("D117").Value = "Bathroom"
$ws.Range("E117").Value = "22.9C"
$ws.Range("F117").Value = "Active"
$ws.Range("A118").Value = "2026-01-28"
$ws.Range("B118").Value = "16:46:50"
$ws.Range("C118").Value = "16:00"
$ws.Range("D118").Value = "Bathroom"
$ws.Range("E118").Value = "22.8C"
$ws.Range("F118").Value = "Active"
$ws.Range("A119").Value = "2026-01-28"
$ws.Range("B119").Value = "16:46:53"
$ws.Range("C119").Value = "16:00"
$ws.Range("D119").Value = "Bathroom"
$ws.Range("E119").Value = "22.8C"
$ws.Range("F119").Value = "Active"
$ws.Range("A120").Value = "2026-01-28"
$ws.Range("B120").Value = "16:47:01"
$ws.Range("C120").Value = "16:00"
$ws.Range("D120").Value = "Bathroom"
$ws.Range("E120").Value = "22.8C"
$ws.Range("F120").Value = "Active"
$ws.Range("A121").Value = "2026-01-28"
$ws.Range("B121").Value = "16:47:05"
$ws.Range("C121").Value = "16:00"
$ws.Range("D121").Value = "Bathroom"
$ws.Range("E121").Value = "22.8C"
$ws.Range("F121").Value = "Active"
$ws.Range("A122").Value = "2026-01-28"
$ws.Range("B122").Value = "16:47:09"
$ws.Range("C122").Value = "16:00"
$ws.Range("D122").Value = "Bathroom"
$ws.Range("E122").Value = "22.8C"
$ws.Range("F122").Value = "Active"
$ws.Range("A123").Value = "2026-01-28"
$ws.Range("B123").Value = "16:47:13"
$ws.Range("C123").Value = "16:00"
$ws.Range("D123").Value = "Bathroom"
$ws.Range("E123").Value = "22.9C"
$ws.Range("F123").Value = "Active"
$ws.Range("A124").Value = "2026-01-28"
$ws.Range("B124").Value = "16:47:29"
$ws.Range("C124").Value = "16:00"
$ws.Range("D124").Value = "Bathroom"
$ws.Range("E124").Value = "22.9C"
$ws.Range("F124").Value = "Active"
$ws.Range("A125").Value = "2026-01-28"
$ws.Range("B125").Value = "16:47:34"
$ws.Range("C125").Value = "16:00"
$ws.Range("D125").Value = "Bathroom"
$ws.Range("E125").Value = "22.9C"
$ws.Range("F125").Value = "Active"
$ws.Range("A126").Value = "2026-01-28"
$ws.Range("B126").Value = "16:47:37"
$ws.Range("C126").Value = "16:00"
$ws.Range("D126").Value = "Bathroom"
$ws.Range("E126").Value = "22.9C"
$ws.Range("F126").Value = "Active"
$ws.Range("A127").Value = "2026-01-28"
$ws.Range("B127").Value = "16:47:41"
$ws.Range("C127").Value = "16:00"
$ws.Range("D127").Value = "Bathroom"
$ws.Range("E127").Value = "22.9C"
$ws.Range("F127").Value = "Active"
